# Apply the "Improved Power Point Presentation" edit:
#   1. Footer date auto-field on the slide master + every slide layout
#      flips from 4/8/2020 to 4/9/2020.
#   2. Slide 3 ("FEATURES AND LABELS"): fix the "InThroat" typo so it
#      reads "in Throat" (split across two runs, matching how PowerPoint
#      records an in-place retype).
#   3. Slide 5 ("Future Plan"): two paragraphs that used to be split
#      across two runs (with identical formatting) are retyped as a
#      single contiguous run each.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder (ppPlaceholderDate = 16) on master + all layouts
# ---------------------------------------------------------------------
function Find-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $sh = $shapes.Item($k)
        try {
            $pf = $sh.PlaceholderFormat
            if ($pf.Type -eq 16) {
                return $sh
            }
        } catch {
            # not a placeholder at all - skip
        }
    }
    return $null
}

$master = $p.SlideMaster

$masterDateShape = Find-DatePlaceholder $master.Shapes
if ($masterDateShape -ne $null) {
    $masterDateShape.TextFrame.TextRange.Text = "4/9/2020"
}

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $dateShape = Find-DatePlaceholder $layout.Shapes
    if ($dateShape -ne $null) {
        $dateShape.TextFrame.TextRange.Text = "4/9/2020"
    }
}

# ---------------------------------------------------------------------
# 2. Slide 3: "Pain InThroat" -> "Pain in Throat"
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(1)
$tr3 = $shape3.TextFrame.TextRange

$para3 = $tr3.Paragraphs(8, 1)
$run3b = $para3.Runs(2, 1)
$runStart = $run3b.Start

# Replace the leading "I" with a lower-case "i" ...
$firstChar = $tr3.Characters($runStart, 1)
$firstChar.Text = "i"

# ... then fix up the remainder ("nThroat" -> "n Throat").
$tr3b = $shape3.TextFrame.TextRange
$rest = $tr3b.Characters($runStart + 1, 7)
$rest.Text = "n Throat"

# ---------------------------------------------------------------------
# 3. Slide 5: merge runs that share formatting into a single run
# ---------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$shape5 = $slide5.Shapes.Item(1)
$tr5 = $shape5.TextFrame.TextRange

# Paragraph 1: "...current situation" + "." -> "...current situation."
$para1 = $tr5.Paragraphs(1, 1)
$p1run2 = $para1.Runs(2, 1)
$p1run2Start = $p1run2.Start
$p1run1Chars = $tr5.Characters($para1.Start, $p1run2Start - $para1.Start)
$p1run1Chars.Text = "Many other ml engineers can work on this model. I have just written very few lines of code but this idea can be used and it has the potential to improve the current situation."

$tr5b = $shape5.TextFrame.TextRange
$para1b = $tr5b.Paragraphs(1, 1)
$p1run2b = $para1b.Runs(2, 1)
$p1run2b.Text = ""

# Paragraph 3: "By fine tu" + "ning the model ... we " -> one run
$tr5c = $shape5.TextFrame.TextRange
$para3b = $tr5c.Paragraphs(3, 1)
$p3run1 = $para3b.Runs(1, 1)
$p3run2 = $para3b.Runs(2, 1)
$p3run1Start = $p3run1.Start
$p3run2Start = $p3run2.Start
$p3run1Chars = $tr5c.Characters($p3run1Start, $p3run2Start - $p3run1Start)
$p3run1Chars.Text = "By fine tuning the model and by trying out different classification algorithms, we "

$tr5d = $shape5.TextFrame.TextRange
$para3c = $tr5d.Paragraphs(3, 1)
$p3run2b = $para3c.Runs(2, 1)
$p3run2b.Text = ""
